$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "fantasy points" column (E) values before touching
# anything, then re-write them into the new column G further down. We
# deliberately avoid Insert/Cut (which re-serializes already-stored floats
# through a lossier formatter) and instead just assign fresh .Value's to
# every affected cell.
$fantasyPoints = @{}
for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints[$r] = $ws.Cells.Item($r, 5).Value2
}

# Headers: E=height, F=weight, G=fantasy points (moved from old E)
$ws.Cells.Item(1, 5).Value = "height"
$ws.Cells.Item(1, 6).Value = "weight"
$ws.Cells.Item(1, 7).Value = "fantasy points"

# Match the bold/bordered header style used by the other header cells
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 260
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}
